$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 1
$ws.Range("F3").Value = 2
$ws.Range("F4").Value = 0
$ws.Range("F5").Value = 0
$ws.Range("F6").Value = 2
$ws.Range("F7").Value = 1
$ws.Range("F8").Value = 1
$ws.Range("F9").Value = 1
$ws.Range("F10").Value = 2
$ws.Range("F11").Value = 0
$ws.Range("F12").Value = 2
$ws.Range("F13").Value = 1
$ws.Range("F14").Value = 2
$ws.Range("F15").Value = 2
$ws.Range("F16").Value = 0
$ws.Range("F17").Value = 2
$ws.Range("F18").Value = 2
$ws.Range("F19").Value = 1
$ws.Range("F20").Value = 0
$ws.Range("F21").Value = 0
$ws.Range("F22").Value = 1
$ws.Range("F23").Value = 0
$ws.Range("F24").Value = 0
$ws.Range("F25").Value = 1
$ws.Range("F26").Value = 1
$ws.Range("F27").Value = 1
$ws.Range("F28").Value = 2
$ws.Range("F29").Value = 1
$ws.Range("F30").Value = 1
$ws.Range("F31").Value = 2
$ws.Range("F32").Value = 0
$ws.Range("F33").Value = 2
$ws.Range("F34").Value = 1
$ws.Range("F35").Value = 1
$ws.Range("F36").Value = 2
$ws.Range("F37").Value = 0
$ws.Range("F38").Value = 0
$ws.Range("F39").Value = 1
$ws.Range("F40").Value = 1
$ws.Range("F41").Value = 0
$ws.Range("F42").Value = 2
$ws.Range("F43").Value = 2
$ws.Range("F45").Value = 2
$ws.Range("F46").Value = 1
$ws.Range("F47").Value = 1
$ws.Range("F48").Value = 1
$ws.Range("F49").Value = 1
$ws.Range("F50").Value = 1
$ws.Range("F51").Value = 1
$ws.Range("F52").Value = 2
$ws.Range("F53").Value = 1
$ws.Range("F54").Value = 1
$ws.Range("F55").Value = 2
$ws.Range("F56").Value = 0
$ws.Range("F57").Value = 2
$ws.Range("F58").Value = 1
$ws.Range("F59").Value = 2
$ws.Range("F61").Value = 2
$ws.Range("F62").Value = 0
$ws.Range("F63").Value = 2
$ws.Range("F64").Value = 2
$ws.Range("F65").Value = 2
$ws.Range("F66").Value = 0
$ws.Range("F67").Value = 2
$ws.Range("F68").Value = 0
$ws.Range("F69").Value = 0
$ws.Range("F70").Value = 2
$ws.Range("F71").Value = 2
$ws.Range("F72").Value = 2
$ws.Range("F73").Value = 1
$ws.Range("F74").Value = 2
$ws.Range("F75").Value = 2
$ws.Range("F76").Value = 2
$ws.Range("F77").Value = 2
$ws.Range("F78").Value = 0
$ws.Range("F79").Value = 1
$ws.Range("F80").Value = 1
$ws.Range("F81").Value = 2
$ws.Range("F82").Value = 0
$ws.Range("F83").Value = 0
$ws.Range("F84").Value = 1
$ws.Range("F85").Value = 2
$ws.Range("F86").Value = 2
$ws.Range("F87").Value = 1
$ws.Range("F88").Value = 2
$ws.Range("F89").Value = 1
$ws.Range("F90").Value = 0
$ws.Range("F91").Value = 1
$ws.Range("F92").Value = 1
$ws.Range("F93").Value = 1
$ws.Range("F94").Value = 1
$ws.Range("F95").Value = 0
$ws.Range("F96").Value = 1
$ws.Range("F97").Value = 2
$ws.Range("F98").Value = 1
$ws.Range("F99").Value = 0
$ws.Range("F100").Value = 0
$ws.Range("F101").Value = 1
$ws.Range("F102").Value = 1
$ws.Range("F103").Value = 0
$ws.Range("F104").Value = 2
$ws.Range("F105").Value = 1
$ws.Range("F106").Value = 0
$ws.Range("F107").Value = 1
$ws.Range("F108").Value = 2
$ws.Range("F109").Value = 2
$ws.Range("F110").Value = 0
$ws.Range("F111").Value = 1
$ws.Range("F112").Value = 0
$ws.Range("F113").Value = 1
$ws.Range("F114").Value = 0
$ws.Range("F115").Value = 1
$ws.Range("F116").Value = 1
$ws.Range("F117").Value = 1
$ws.Range("F118").Value = 0
$ws.Range("F119").Value = 0
$ws.Range("F120").Value = 0
$ws.Range("F121").Value = 0
$ws.Range("F122").Value = 2
$ws.Range("F123").Value = 1
$ws.Range("F124").Value = 2
$ws.Range("F125").Value = 1
$ws.Range("F126").Value = 1
$ws.Range("F127").Value = 1
$ws.Range("F128").Value = 1
$ws.Range("F129").Value = 0
$ws.Range("F130").Value = 0
$ws.Range("F131").Value = 2
$ws.Range("F132").Value = 2
$ws.Range("F133").Value = 0
$ws.Range("F134").Value = 0
$ws.Range("F135").Value = 1
$ws.Range("F136").Value = 1
$ws.Range("F137").Value = 0
$ws.Range("F138").Value = 2
$ws.Range("F139").Value = 2
$ws.Range("F140").Value = 2
$ws.Range("F141").Value = 2
$ws.Range("F142").Value = 2
$ws.Range("F143").Value = 0
$ws.Range("F144").Value = 0
